$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the MIN @ MEM row (Jan 15 game that was pulled / not used); this shifts
# rows 77-79 up to 76-78.
$ws.Rows(76).Delete()

# Fill in the "Beat Vegas?" results for the Jan 15 slate now that the games
# have been played.
$ws.Range("G72").Value = "No"
$ws.Range("G73").Value = "Yes"
$ws.Range("G74").Value = "No"
$ws.Range("G75").Value = "No"
$ws.Range("G76").Value = "Yes"
$ws.Range("G77").Value = "Yes"
$ws.Range("G78").Value = "No"

# Ran the model for Jan 16 -- append the predictions for that day's games.
$ws.Range("A79").Value = 44212
$ws.Range("B79").Value = "SAS"
$ws.Range("C79").Value = "HOU"
$ws.Range("D79").Value = -7
$ws.Range("E79").Value = 5.6
$ws.Range("F79").Value = -12.6

$ws.Range("A80").Value = 44212
$ws.Range("B80").Value = "BRK"
$ws.Range("C80").Value = "ORL"
$ws.Range("D80").Value = -8
$ws.Range("E80").Value = -20.9
$ws.Range("F80").Value = 12.9

$ws.Range("A81").Value = 44212
$ws.Range("B81").Value = "TOR"
$ws.Range("C81").Value = "CHO"
$ws.Range("D81").Value = -7
$ws.Range("E81").Value = -4
$ws.Range("F81").Value = -3

$ws.Range("A82").Value = 44212
$ws.Range("B82").Value = "MEM"
$ws.Range("C82").Value = "PHI"
$ws.Range("D82").Value = -1.5
$ws.Range("E82").Value = 8.2
$ws.Range("F82").Value = -9.7

$ws.Range("A83").Value = 44212
$ws.Range("B83").Value = "MIA"
$ws.Range("C83").Value = "DET"
$ws.Range("D83").Value = -5
$ws.Range("E83").Value = -18.6
$ws.Range("F83").Value = 13.6

$ws.Range("A84").Value = 44212
$ws.Range("B84").Value = "POR"
$ws.Range("C84").Value = "ATL"
$ws.Range("D84").Value = -5
$ws.Range("E84").Value = -12.4
$ws.Range("F84").Value = 7.4

# The new date cells should carry the same date-formatted style used by the
# rest of column A.
$ws.Range("A79:A84").NumberFormat = $ws.Range("A78").NumberFormat

# Restore the cursor/selection to where the author left off working.
[void]$ws.Range("L93").Select()
